$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("sort")
$ws2.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws2.Range("A12").Select()
